$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "51.723.59"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "'" + "3.108.98"
$ws.Range("E3").Value = "  +3.92%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'" + "388.59"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").Value = "'" + "103.29"
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("D7").Value = "'" + "0.543"
$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'" + "0.587"
$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("D10").Value = "'" + "36.97"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "'" + "0.0861"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").Value = "'" + "3.592.03"
$ws.Range("E13").Value = "  +3.90%  "

$ws.Range("D14").Value = "'" + "18.67"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").Value = "'" + "7.87"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").Value = "'" + "3.108.85"
$ws.Range("E16").Value = "  +4.01%  "

$ws.Range("D17").Value = "'" + "0.982"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").Value = "'" + "10.74"
$ws.Range("E18").Value = "  -4.10%  "

$ws.Range("D19").Value = "'" + "51.848.69"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "'" + "3.18"
$ws.Range("E20").Value = "  +2.68%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'" + "0.0" + ([string][char]8323) + "0968"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'" + "12.42"
$ws.Range("E22").Value = "  -1.47%  "

$ws.Range("D23").Value = "'" + "70.05"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").Value = "'" + "268.13"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  -3.63%  "

$ws.Range("D26").Value = "'" + "8.07"
$ws.Range("E26").Value = "  +1.84%  "

$ws.Range("D27").Value = "'" + "27.05"
$ws.Range("E27").Value = "  +3.44%  "

$ws.Range("D28").Value = "'" + "0.170"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").Value = "'" + "7.20"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("D32").Value = "'" + "10.33"
$ws.Range("E32").Value = "  -1.06%  "

$ws.Range("D33").Value = "'" + "35.38"
$ws.Range("E33").Value = "  +1.93%  "

$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'" + "50.41"
$ws.Range("E35").Value = "  -2.00%  "

$ws.Range("D36").Value = "'" + "0.0449"
$ws.Range("E36").Value = "  +0.50%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "'" + "3.43"
$ws.Range("E38").Value = "  +4.06%  "

$ws.Range("D39").Value = "'" + "0.289"
$ws.Range("E39").Value = "  +6.26%  "

$ws.Range("E40").Value = "  +1.82%  "

$ws.Range("D41").Value = "'" + "16.88"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").Value = "'" + "128.07"
$ws.Range("E43").Value = "  +1.17%  "

$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("D45").Value = "'" + "3.68"
$ws.Range("E45").Value = "  -4.93%  "

$ws.Range("D46").Value = "'" + "22.28"
$ws.Range("E46").Value = "  +3.97%  "

$ws.Range("D47").Value = "'" + "2.50"
$ws.Range("E47").Value = "  +6.28%  "

$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("D49").Value = "'" + "2.046.95"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").Value = "'" + "3.416.14"
$ws.Range("E50").Value = "  +4.06%  "

$ws.Range("D51").Value = "'" + "0.206"
$ws.Range("E51").Value = "  +5.24%  "
